$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 37.17328633333333
$ws.Cells.Item(2, 8).Value = 111.519859
$ws.Cells.Item(2, 9).Value = 0.005170079968594893
$ws.Cells.Item(2, 10).Value = 0.005188590814393131
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.6435283333333334
$ws.Cells.Item(2, 14).Value = 1.930585
$ws.Cells.Item(2, 15).Value = 0.1384187503011309
$ws.Cells.Item(2, 16).Value = 0.1430362385488992
$ws.Cells.Item(2, 17).Value = 23.92206299861278
$ws.Cells.Item(2, 18).Value = 215.298566987515
$ws.Cells.Item(2, 19).Value = 0.0007156360082098151
$ws.Cells.Item(2, 20).Value = 0.000742156513460163

$ws.Cells.Item(3, 7).Value = 37.17328633333333
$ws.Cells.Item(3, 8).Value = 111.519859
$ws.Cells.Item(3, 9).Value = 0.005170079968594893
$ws.Cells.Item(3, 10).Value = 0.005188590814393131
$ws.Cells.Item(3, 15).Value = 0.2324275794584003
$ws.Cells.Item(3, 16).Value = 0.2401810927235583
$ws.Cells.Item(3, 17).Value = 40.16903191455489
$ws.Cells.Item(3, 18).Value = 361.521287230994
$ws.Cells.Item(3, 19).Value = 0.001201669172706873
$ws.Cells.Item(3, 20).Value = 0.001246201411496359

$ws.Cells.Item(4, 7).Value = 37.17328633333333
$ws.Cells.Item(4, 8).Value = 111.519859
$ws.Cells.Item(4, 9).Value = 0.005170079968594893
$ws.Cells.Item(4, 10).Value = 0.005188590814393131
$ws.Cells.Item(4, 13).Value = 1.770781666666667
$ws.Cells.Item(4, 14).Value = 5.312345000000001
$ws.Cells.Item(4, 15).Value = 0.3808835954223518
$ws.Cells.Item(4, 16).Value = 0.3935894284240538
$ws.Cells.Item(4, 17).Value = 65.82577392881723
$ws.Cells.Item(4, 18).Value = 592.431965359355
$ws.Cells.Item(4, 19).Value = 0.001969198647059502
$ws.Cells.Item(4, 20).Value = 0.002042174492963288

$ws.Cells.Item(5, 7).Value = 37.17328633333333
$ws.Cells.Item(5, 8).Value = 111.519859
$ws.Cells.Item(5, 9).Value = 0.005170079968594893
$ws.Cells.Item(5, 10).Value = 0.005188590814393131
$ws.Cells.Item(5, 13).Value = 0.45025
$ws.Cells.Item(5, 14).Value = 0.9005000000000001
$ws.Cells.Item(5, 15).Value = 0.09684584049355638
$ws.Cells.Item(5, 16).Value = 0.06671766993594362
$ws.Cells.Item(5, 17).Value = 16.73727217158333
$ws.Cells.Item(5, 18).Value = 100.4236330295
$ws.Cells.Item(5, 19).Value = 0.000500700739977472
$ws.Cells.Item(5, 20).Value = 0.0003461706893873498

$ws.Cells.Item(6, 7).Value = 37.17328633333333
$ws.Cells.Item(6, 8).Value = 111.519859
$ws.Cells.Item(6, 9).Value = 0.005170079968594893
$ws.Cells.Item(6, 10).Value = 0.005188590814393131
$ws.Cells.Item(6, 13).Value = 0.7039926666666667
$ws.Cells.Item(6, 14).Value = 2.111978
$ws.Cells.Item(6, 15).Value = 0.1514242343245606
$ws.Cells.Item(6, 16).Value = 0.1564755703675451
$ws.Cells.Item(6, 17).Value = 26.16972097456689
$ws.Cells.Item(6, 18).Value = 235.527488771102
$ws.Cells.Item(6, 19).Value = 0.0007828754006412299
$ws.Cells.Item(6, 20).Value = 0.0008118877070859705

$ws.Cells.Item(7, 9).Value = 0.006280726092526873
$ws.Cells.Item(7, 10).Value = 0.006303213472394487
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.6435283333333334
$ws.Cells.Item(7, 14).Value = 1.930585
$ws.Cells.Item(7, 15).Value = 0.1384187503011309
$ws.Cells.Item(7, 16).Value = 0.1430362385488992
$ws.Cells.Item(7, 17).Value = 29.06104473724278
$ws.Cells.Item(7, 18).Value = 261.549402635185
$ws.Cells.Item(7, 19).Value = 0.0008693702567112746
$ws.Cells.Item(7, 20).Value = 0.0009015879458620531

$ws.Cells.Item(8, 9).Value = 0.006280726092526873
$ws.Cells.Item(8, 10).Value = 0.006303213472394487
$ws.Cells.Item(8, 15).Value = 0.2324275794584003
$ws.Cells.Item(8, 16).Value = 0.2401810927235583
$ws.Cells.Item(8, 19).Value = 0.001459813962927238
$ws.Cells.Item(8, 20).Value = 0.001513912699469562

$ws.Cells.Item(9, 9).Value = 0.006280726092526873
$ws.Cells.Item(9, 10).Value = 0.006303213472394487
$ws.Cells.Item(9, 13).Value = 1.770781666666667
$ws.Cells.Item(9, 14).Value = 5.312345000000001
$ws.Cells.Item(9, 15).Value = 0.3808835954223518
$ws.Cells.Item(9, 16).Value = 0.3935894284240538
$ws.Cells.Item(9, 17).Value = 79.96658821272723
$ws.Cells.Item(9, 18).Value = 719.6992939145451
$ws.Cells.Item(9, 19).Value = 0.002392225535984614
$ws.Cells.Item(9, 20).Value = 0.002480878187834542

$ws.Cells.Item(10, 9).Value = 0.006280726092526873
$ws.Cells.Item(10, 10).Value = 0.006303213472394487
$ws.Cells.Item(10, 13).Value = 0.45025
$ws.Cells.Item(10, 14).Value = 0.9005000000000001
$ws.Cells.Item(10, 15).Value = 0.09684584049355638
$ws.Cells.Item(10, 16).Value = 0.06671766993594362
$ws.Cells.Item(10, 17).Value = 20.33280388008334
$ws.Cells.Item(10, 18).Value = 121.9968232805
$ws.Cells.Item(10, 19).Value = 0.0006082621973405751
$ws.Cells.Item(10, 20).Value = 0.0004205357159870085

$ws.Cells.Item(11, 9).Value = 0.006280726092526873
$ws.Cells.Item(11, 10).Value = 0.006303213472394487
$ws.Cells.Item(11, 13).Value = 0.7039926666666667
$ws.Cells.Item(11, 14).Value = 2.111978
$ws.Cells.Item(11, 15).Value = 0.1514242343245606
$ws.Cells.Item(11, 16).Value = 0.1564755703675451
$ws.Cells.Item(11, 17).Value = 31.79154874925089
$ws.Cells.Item(11, 18).Value = 286.1239387432581
$ws.Cells.Item(11, 19).Value = 0.000951054139563171
$ws.Cells.Item(11, 20).Value = 0.0009862989232413218

$ws.Cells.Item(12, 7).Value = 2375.59786
$ws.Cells.Item(12, 8).Value = 7126.79358
$ws.Cells.Item(12, 9).Value = 0.3303993840977568
$ws.Cells.Item(12, 10).Value = 0.3315823391174117
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.6435283333333334
$ws.Cells.Item(12, 14).Value = 1.930585
$ws.Cells.Item(12, 15).Value = 0.1384187503011309
$ws.Cells.Item(12, 16).Value = 0.1430362385488992
$ws.Cells.Item(12, 17).Value = 1528.764531516033
$ws.Cells.Item(12, 18).Value = 13758.8807836443
$ws.Cells.Item(12, 19).Value = 0.04573346984707483
$ws.Cells.Item(12, 20).Value = 0.04742829055660008

$ws.Cells.Item(13, 7).Value = 2375.59786
$ws.Cells.Item(13, 8).Value = 7126.79358
$ws.Cells.Item(13, 9).Value = 0.3303993840977568
$ws.Cells.Item(13, 10).Value = 0.3315823391174117
$ws.Cells.Item(13, 15).Value = 0.2324275794584003
$ws.Cells.Item(13, 16).Value = 0.2401810927235583
$ws.Cells.Item(13, 17).Value = 2567.044124073586
$ws.Cells.Item(13, 18).Value = 23103.39711666228
$ws.Cells.Item(13, 19).Value = 0.0767939291003879
$ws.Cells.Item(13, 20).Value = 0.07963980853705341

$ws.Cells.Item(14, 7).Value = 2375.59786
$ws.Cells.Item(14, 8).Value = 7126.79358
$ws.Cells.Item(14, 9).Value = 0.3303993840977568
$ws.Cells.Item(14, 10).Value = 0.3315823391174117
$ws.Cells.Item(14, 13).Value = 1.770781666666667
$ws.Cells.Item(14, 14).Value = 5.312345000000001
$ws.Cells.Item(14, 15).Value = 0.3808835954223518
$ws.Cells.Item(14, 16).Value = 0.3935894284240538
$ws.Cells.Item(14, 17).Value = 4206.665137860567
$ws.Cells.Item(14, 18).Value = 37859.9862407451
$ws.Cells.Item(14, 19).Value = 0.1258437053404842
$ws.Cells.Item(14, 20).Value = 0.1305073033287328

$ws.Cells.Item(15, 7).Value = 2375.59786
$ws.Cells.Item(15, 8).Value = 7126.79358
$ws.Cells.Item(15, 9).Value = 0.3303993840977568
$ws.Cells.Item(15, 10).Value = 0.3315823391174117
$ws.Cells.Item(15, 13).Value = 0.45025
$ws.Cells.Item(15, 14).Value = 0.9005000000000001
$ws.Cells.Item(15, 15).Value = 0.09684584049355638
$ws.Cells.Item(15, 16).Value = 0.06671766993594362
$ws.Cells.Item(15, 17).Value = 1069.612936465
$ws.Cells.Item(15, 18).Value = 6417.67761879
$ws.Cells.Item(15, 19).Value = 0.03199780605150063
$ws.Cells.Item(15, 20).Value = 0.0221224010578236

$ws.Cells.Item(16, 7).Value = 2375.59786
$ws.Cells.Item(16, 8).Value = 7126.79358
$ws.Cells.Item(16, 9).Value = 0.3303993840977568
$ws.Cells.Item(16, 10).Value = 0.3315823391174117
$ws.Cells.Item(16, 13).Value = 0.7039926666666667
$ws.Cells.Item(16, 14).Value = 2.111978
$ws.Cells.Item(16, 15).Value = 0.1514242343245606
$ws.Cells.Item(16, 16).Value = 0.1564755703675451
$ws.Cells.Item(16, 17).Value = 1672.403472389027
$ws.Cells.Item(16, 18).Value = 15051.63125150124
$ws.Cells.Item(16, 19).Value = 0.05003047375830923
$ws.Cells.Item(16, 20).Value = 0.05188453563720175

$ws.Cells.Item(17, 7).Value = 76.954105
$ws.Cells.Item(17, 8).Value = 153.90821
$ws.Cells.Item(17, 9).Value = 0.01070281688829022
$ws.Cells.Item(17, 10).Value = 0.007160758019481436
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.6435283333333334
$ws.Cells.Item(17, 14).Value = 1.930585
$ws.Cells.Item(17, 15).Value = 0.1384187503011309
$ws.Cells.Item(17, 16).Value = 0.1430362385488992
$ws.Cells.Item(17, 17).Value = 49.52214693380834
$ws.Cells.Item(17, 18).Value = 297.13288160285
$ws.Cells.Item(17, 19).Value = 0.001481470538378971
$ws.Cells.Item(17, 20).Value = 0.00102424789226549

$ws.Cells.Item(18, 7).Value = 76.954105
$ws.Cells.Item(18, 8).Value = 153.90821
$ws.Cells.Item(18, 9).Value = 0.01070281688829022
$ws.Cells.Item(18, 10).Value = 0.007160758019481436
$ws.Cells.Item(18, 15).Value = 0.2324275794584003
$ws.Cells.Item(18, 16).Value = 0.2401810927235583
$ws.Cells.Item(18, 17).Value = 83.15573371647666
$ws.Cells.Item(18, 18).Value = 498.93440229886
$ws.Cells.Item(18, 19).Value = 0.002487629822731784
$ws.Cells.Item(18, 20).Value = 0.001719878685848034

$ws.Cells.Item(19, 7).Value = 76.954105
$ws.Cells.Item(19, 8).Value = 153.90821
$ws.Cells.Item(19, 9).Value = 0.01070281688829022
$ws.Cells.Item(19, 10).Value = 0.007160758019481436
$ws.Cells.Item(19, 13).Value = 1.770781666666667
$ws.Cells.Item(19, 14).Value = 5.312345000000001
$ws.Cells.Item(19, 15).Value = 0.3808835954223518
$ws.Cells.Item(19, 16).Value = 0.3935894284240538
$ws.Cells.Item(19, 17).Value = 136.2689183087417
$ws.Cells.Item(19, 18).Value = 817.6135098524501
$ws.Cells.Item(19, 19).Value = 0.004076527377559047
$ws.Cells.Item(19, 20).Value = 0.002818398655970658

$ws.Cells.Item(20, 7).Value = 76.954105
$ws.Cells.Item(20, 8).Value = 153.90821
$ws.Cells.Item(20, 9).Value = 0.01070281688829022
$ws.Cells.Item(20, 10).Value = 0.007160758019481436
$ws.Cells.Item(20, 13).Value = 0.45025
$ws.Cells.Item(20, 14).Value = 0.9005000000000001
$ws.Cells.Item(20, 15).Value = 0.09684584049355638
$ws.Cells.Item(20, 16).Value = 0.06671766993594362
$ws.Cells.Item(20, 17).Value = 34.64858577625
$ws.Cells.Item(20, 18).Value = 138.594343105
$ws.Cells.Item(20, 19).Value = 0.001036523297195096
$ws.Cells.Item(20, 20).Value = 0.0004777490900349238

$ws.Cells.Item(21, 7).Value = 76.954105
$ws.Cells.Item(21, 8).Value = 153.90821
$ws.Cells.Item(21, 9).Value = 0.01070281688829022
$ws.Cells.Item(21, 10).Value = 0.007160758019481436
$ws.Cells.Item(21, 13).Value = 0.7039926666666667
$ws.Cells.Item(21, 14).Value = 2.111978
$ws.Cells.Item(21, 15).Value = 0.1514242343245606
$ws.Cells.Item(21, 16).Value = 0.1564755703675451
$ws.Cells.Item(21, 17).Value = 54.17512558989667
$ws.Cells.Item(21, 18).Value = 325.05075353938
$ws.Cells.Item(21, 19).Value = 0.001620665852425323
$ws.Cells.Item(21, 20).Value = 0.00112048369536233

$ws.Cells.Item(22, 7).Value = 4655.195393666666
$ws.Cells.Item(22, 8).Value = 13965.586181
$ws.Cells.Item(22, 9).Value = 0.6474469929528313
$ws.Cells.Item(22, 10).Value = 0.6497650985763194
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 0.6435283333333334
$ws.Cells.Item(22, 14).Value = 1.930585
$ws.Cells.Item(22, 15).Value = 0.1384187503011309
$ws.Cells.Item(22, 16).Value = 0.1430362385488992
$ws.Cells.Item(22, 17).Value = 2995.75013302732
$ws.Cells.Item(22, 18).Value = 26961.75119724588
$ws.Cells.Item(22, 19).Value = 0.08961880365075599
$ws.Cells.Item(22, 20).Value = 0.09293995564071141

$ws.Cells.Item(23, 7).Value = 4655.195393666666
$ws.Cells.Item(23, 8).Value = 13965.586181
$ws.Cells.Item(23, 9).Value = 0.6474469929528313
$ws.Cells.Item(23, 10).Value = 0.6497650985763194
$ws.Cells.Item(23, 15).Value = 0.2324275794584003
$ws.Cells.Item(23, 16).Value = 0.2401810927235583
$ws.Cells.Item(23, 17).Value = 5030.351383515071
$ws.Cells.Item(23, 18).Value = 45273.16245163565
$ws.Cells.Item(23, 19).Value = 0.1504845373996465
$ws.Cells.Item(23, 20).Value = 0.1560612913896909

$ws.Cells.Item(24, 7).Value = 4655.195393666666
$ws.Cells.Item(24, 8).Value = 13965.586181
$ws.Cells.Item(24, 9).Value = 0.6474469929528313
$ws.Cells.Item(24, 10).Value = 0.6497650985763194
$ws.Cells.Item(24, 13).Value = 1.770781666666667
$ws.Cells.Item(24, 14).Value = 5.312345000000001
$ws.Cells.Item(24, 15).Value = 0.3808835954223518
$ws.Cells.Item(24, 16).Value = 0.3935894284240538
$ws.Cells.Item(24, 17).Value = 8243.334657856049
$ws.Cells.Item(24, 18).Value = 74190.01192070445
$ws.Cells.Item(24, 19).Value = 0.2466019385212645
$ws.Cells.Item(24, 20).Value = 0.2557406737585525

$ws.Cells.Item(25, 7).Value = 4655.195393666666
$ws.Cells.Item(25, 8).Value = 13965.586181
$ws.Cells.Item(25, 9).Value = 0.6474469929528313
$ws.Cells.Item(25, 10).Value = 0.6497650985763194
$ws.Cells.Item(25, 13).Value = 0.45025
$ws.Cells.Item(25, 14).Value = 0.9005000000000001
$ws.Cells.Item(25, 15).Value = 0.09684584049355638
$ws.Cells.Item(25, 16).Value = 0.06671766993594362
$ws.Cells.Item(25, 17).Value = 2096.001725998417
$ws.Cells.Item(25, 18).Value = 12576.0103559905
$ws.Cells.Item(25, 19).Value = 0.06270254820754262
$ws.Cells.Item(25, 20).Value = 0.04335081338271075

$ws.Cells.Item(26, 7).Value = 4655.195393666666
$ws.Cells.Item(26, 8).Value = 13965.586181
$ws.Cells.Item(26, 9).Value = 0.6474469929528313
$ws.Cells.Item(26, 10).Value = 0.6497650985763194
$ws.Cells.Item(26, 13).Value = 0.7039926666666667
$ws.Cells.Item(26, 14).Value = 2.111978
$ws.Cells.Item(26, 15).Value = 0.1514242343245606
$ws.Cells.Item(26, 16).Value = 0.1564755703675451
$ws.Cells.Item(26, 17).Value = 3277.223419041779
$ws.Cells.Item(26, 18).Value = 29495.01077137602
$ws.Cells.Item(26, 19).Value = 0.09803916517362166
$ws.Cells.Item(26, 20).Value = 0.1016723644046537
